# Insert two new weekly records (dated 2022-05-13 / serial 44694) at the top
# of the "Pintón" / "Primera Pintón" price series for
# Macroferia Regional de Talca - Plátano, pushing the existing history
# down by two rows (dimension grows from A1:T645 to A1:T647).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows above the current row 544 (Excel shifts everything
# below down, and copies formatting - e.g. the date-number-format on
# column D - from the row above, same as a manual row insert).
$ws.Rows("544:545").Insert()

# --- New row 544: "Pintón" ---
$ws.Cells.Item(544, 1).Value = 5
$ws.Cells.Item(544, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(544, 3).Value = "Maule"
$ws.Cells.Item(544, 4).Value = 44694
$ws.Cells.Item(544, 5).Value = 7
$ws.Cells.Item(544, 6).Value = "Fruta"
$ws.Cells.Item(544, 7).Value = 100108
$ws.Cells.Item(544, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(544, 9).Value = 100108006
$ws.Cells.Item(544, 10).Value = "Plátano"
$ws.Cells.Item(544, 11).Value = "Sin especificar"
$ws.Cells.Item(544, 12).Value = "Pintón"
$ws.Cells.Item(544, 13).Value = 800
$ws.Cells.Item(544, 14).Value = 10000
$ws.Cells.Item(544, 15).Value = 10000
$ws.Cells.Item(544, 16).Value = 10000
$ws.Cells.Item(544, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(544, 18).Value = "Ecuador"
$ws.Cells.Item(544, 19).Value = 500
$ws.Cells.Item(544, 20).Value = 20

# --- New row 545: "Primera Pintón" ---
$ws.Cells.Item(545, 1).Value = 5
$ws.Cells.Item(545, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(545, 3).Value = "Maule"
$ws.Cells.Item(545, 4).Value = 44694
$ws.Cells.Item(545, 5).Value = 7
$ws.Cells.Item(545, 6).Value = "Fruta"
$ws.Cells.Item(545, 7).Value = 100108
$ws.Cells.Item(545, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(545, 9).Value = 100108006
$ws.Cells.Item(545, 10).Value = "Plátano"
$ws.Cells.Item(545, 11).Value = "Sin especificar"
$ws.Cells.Item(545, 12).Value = "Primera Pintón"
$ws.Cells.Item(545, 13).Value = 450
$ws.Cells.Item(545, 14).Value = 11000
$ws.Cells.Item(545, 15).Value = 11000
$ws.Cells.Item(545, 16).Value = 11000
$ws.Cells.Item(545, 17).Value = '$/caja 20 kilos'
$ws.Cells.Item(545, 18).Value = "Ecuador"
$ws.Cells.Item(545, 19).Value = 550
$ws.Cells.Item(545, 20).Value = 20
